$d = $word.ActiveDocument

$replacements = @(
    @{
        Old = 'A imagem apresenta a tela inicial da Google Play Store, mostrando uma seção personalizada chamada "Para você". Ela lista sugestões de jogos baseadas na atividade recente do usuário, incluindo títulos como "Blood Strike", "Roblox" e "Tile Club". Abaixo, há uma seção de jogos patrocinados, com sugestões como "Paciência - Jogo de Solitário", "Coin Master" e "Bubble Pop! Cannon Shooter". A parte inferior da tela exibe ícones para navegar entre Jogos, Apps e Livros, com uma barra de pesquisa destacada em vermelho.'
        New = 'A imagem mostra a interface da Google Play Store, na seção "Para você". No topo, estão destacados jogos como "Blood Strike", "Roblox" e "Tile Club". Abaixo, há uma lista de sugestões de jogos patrocinados, incluindo "Paciência", "Coin Master" e "Bubble Pop! Cannon Shooter", cada um acompanhado de informações de classificação e tamanho. Na parte inferior da tela, há ícones para navegar entre seções, incluindo "Jogos", "Apps", "Livros" e uma opção de "Pesquisar" destacada em vermelho. image_rId8.png'
    },
    @{
        Old = 'A imagem mostra uma tela de pesquisa em um aplicativo de loja digital, onde o usuário procura por "smart sales force". Na parte superior, aparecem os resultados patrocinados relacionados a essa busca. O destaque é para o aplicativo "Smart Força de Vendas" da Arpa Sistemas, que possui uma avaliação de 4,3 estrelas, ocupa 14 MB de espaço e tem mais de mil downloads. Outros aplicativos relacionados, como Salesforce e App Sales Force, também estão listados abaixo.'
        New = 'A imagem mostra uma tela de pesquisa na loja de aplicativos, onde o termo "smart sales force" está sendo utilizado. Os resultados incluem vários aplicativos, com destaque para "Smart Força de Vendas" da Arpa Sistemas, que possui uma classificação de 4,3 estrelas e 14 MB de tamanho, além de mais de mil downloads. Outros aplicativos listados incluem Salesforce, App Sales Force +, e Meta Sales Force, com diferentes classificações e tamanhos. A interface apresenta também um botão de instalação para os aplicativos. image_rId9.png'
    },
    @{
        Old = 'A imagem apresenta a página de download do aplicativo "Smart Força de Vendas" na Google Play Store. O aplicativo, desenvolvido pela Arpa Sistemas, possui uma classificação de 4,2 estrelas, com 12 avaliações e um tamanho de 14 MB. A interface exibe capturas de tela do aplicativo em uso, mostrando diferentes funcionalidades e opções. Há uma seção informativa sobre o aplicativo, destacando que "Smart Vendas" agora é chamado de "Smart Força de Vendas", e um botão para instalá-lo. Na parte inferior da tela, há ícones que permitem acessar jogos, outros aplicativos, e livros, além de uma opção de pesquisa.'
        New = 'A imagem apresenta a interface do aplicativo "Smart Força de Vendas", desenvolvido pela Arpa Sistemas. Na parte superior, está o nome do aplicativo junto com a sua classificação de 4,2 estrelas, o número de avaliações (12) e o tamanho do aplicativo (14 MB). Abaixo, há uma chamada para ação para instalar o aplicativo. A imagem exibe também várias capturas de tela do aplicativo, mostrando suas funcionalidades. Há seções como "Sobre este app" e "Segurança dos dados" também apresentadas na parte inferior. Além disso, são visualizados ícones representando diferentes categorias como jogos, apps, e livros. image_rId10.png'
    },
    @{
        Old = 'A imagem mostra a tela de instalação do aplicativo "Smart Força de Vendas" em um dispositivo móvel. Na parte superior, há uma barra de status com o horário e a qualidade do sinal. Abaixo, são apresentadas sugestões de aplicativos patrocinados, como "Nomad: Conta em Dólar e Cartão", "Livelo: juntar e trocar pontos" e "Estoque, Vendas, PDV, Finanças". Na parte inferior, há uma seção com mais aplicativos para testar, incluindo "PictureThis", "Arquivos do Google" e "CamScanner". A interface é simples, com ícones coloridos representando cada aplicativo.'
        New = 'A imagem exibe a tela de instalação do aplicativo "Smart Força de Vendas" em um dispositivo móvel. Acima, há um botão para cancelar ou abrir o aplicativo, além de um aviso indicando que ele é verificado pelo Play Protect. Abaixo, são apresentadas sugestões de aplicativos patrocinados, como "Nomad: Conta em Dólar e Cartão", "Livelo: juntar e trocar pontos" e "Estoque, Vendas, Pdv, Finanças", juntamente com mais opções de aplicativos para testar, incluindo "PictureThis Identificador Planta" e "CamScanner". A parte inferior da tela contém ícones de acesso a jogos, aplicativos, e livros. image_rId11.png'
    }
)

foreach ($rep in $replacements) {
    $found = $false
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($rep.Old + "`r")) {
            $p.Range.Text = $rep.New
            $found = $true
            break
        }
    }
    if (-not $found) {
        Write-Host "NOT FOUND:" $rep.Old.Substring(0, 40)
    }
}
